$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Construct', ['Token Artifact Creature — Construct', 'Defender', '1/1'])"
$ws.Range("A3").Value = "('Dack Fayden Emblem', ['Emblem — Dack', 'Whenever you cast a spell that targets one or more permanents, gain control of those permanents.'])"
$ws.Range("A4").Value = "('Demon', ['Token Creature — Demon', 'Flying', '*/*'])"
$ws.Range("A5").Value = "('Elephant', ['Token Creature — Elephant', '3/3'])"
$ws.Range("A6").Value = "('Ogre', ['Token Creature — Ogre', '4/4'])"
$ws.Range("A7").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A8").Value = "('Squirrel', ['Token Creature — Squirrel', '1/1'])"
$ws.Range("A9").Value = "('Wolf', ['Token Creature — Wolf', '2/2'])"
$ws.Range("A10").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"

$ws.Range("A11:A31").EntireRow.Delete()
